# PlayerPerformance_4520.xlsx update
# 1. Insert a new "Player Info" sheet before "ODI Batting" with player bio data.
# 2. Rename MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting" and "ODI Bowling" sheets,
#    replacing the full scorecard URL values with just the numeric match code.

$wb = $excel.ActiveWorkbook

# --- Step 1: insert "Player Info" sheet before "ODI Batting" ---
$battingRef = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingRef)
$playerInfo.Name = "Player Info"

# re-fetch by name to get a safe, stable reference
$playerInfo = $wb.Worksheets.Item("Player Info")

$headers = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $headers.Length; $c++) {
    $playerInfo.Cells.Item(1, $c).Value = $headers[$c - 1]
}

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$playerInfo.Cells.Item(2, 1).Value = "'4520"
$playerInfo.Cells.Item(2, 2).Value = "Bradley Thomas James Wheal"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Fast"

# --- Step 2: rename MATCH_CARD_LINK -> MATCH_CODE, and shrink URLs to bare match codes ---

# ODI Batting: MATCH_CARD_LINK lives in column D
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"
for ($r = 2; $r -le 14; $r++) {
    $cell = $batting.Cells.Item($r, 4)
    $code = $cell.Text -replace '.*MatchCode=', ''
    $cell.Value = "'" + $code
}

# ODI Bowling: MATCH_CARD_LINK lives in column B
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"
for ($r = 2; $r -le 14; $r++) {
    $cell = $bowling.Cells.Item($r, 2)
    $code = $cell.Text -replace '.*MatchCode=', ''
    $cell.Value = "'" + $code
}
